$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) : update column F (想去人数) on rows 3-14 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 35
$ws1.Range("F4").Value = 1402
$ws1.Range("F5").Value = 321
$ws1.Range("F6").Value = 1035
$ws1.Range("F7").Value = 10751
$ws1.Range("F8").Value = 24
$ws1.Range("F10").Value = 296
$ws1.Range("F12").Value = 717
$ws1.Range("F13").Value = 12072
$ws1.Range("F14").Value = 12531

# --- Sheet "全部类型" (sheet4) : update column F (想去人数) on rows 4-15 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 35
$ws4.Range("F5").Value = 1402
$ws4.Range("F6").Value = 321
$ws4.Range("F7").Value = 1035
$ws4.Range("F8").Value = 10752
$ws4.Range("F9").Value = 24
$ws4.Range("F11").Value = 296
$ws4.Range("F13").Value = 717
$ws4.Range("F14").Value = 12072
$ws4.Range("F15").Value = 12531
